# Add new columns I (I0) and J (IF) with per-row values, matching the
# target diff which extends the sheet from A1:H55 to A1:J55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(1,5,3,8,7,9,6,4,5,5,8,9,5,8,5,7,8,9,8,9,8,11,7,7,4,8,5,8,9,9,6,7,9,8,7,8,9,7,7,8,7,7,5,6,7,7,6,7,7,3,7,4,4,4)
$jValues = @(2,5,3,8,7,9,7,5,6,6,8,9,6,8,5,8,8,9,8,10,8,11,7,7,5,8,6,8,9,9,6,7,9,8,8,8,9,7,7,8,7,7,5,6,8,7,6,7,7,4,7,4,4,4)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
